# Updated cryptos list values (prices / 1h volume %) per upstream scrape refresh.
# All Price-column (D) literals are forced to text via a leading apostrophe so
# they keep matching the sheets t="inlineStr" cell type instead of Excel
# auto-coercing them to numbers (which would also drop formatting, e.g. "73.80" -> 73.8).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'42.410.28"
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = "'2.180.78"
$ws.Range('E3').Value = '  -1.41%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'254.12"
$ws.Range('E5').Value = '  +5.51%  '
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('D7').Value = "'73.80"
$ws.Range('E7').Value = '  -2.13%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -2.37%  '
$ws.Range('D10').Value = "'40.41"
$ws.Range('E10').Value = '  -2.06%  '
$ws.Range('D11').Value = "'0.0916"
$ws.Range('E11').Value = '  -0.30%  '
$ws.Range('E12').Value = '  +0.52%  '
$ws.Range('D13').Value = "'6.78"
$ws.Range('E13').Value = '  -1.10%  '
$ws.Range('D14').Value = "'2.507.38"
$ws.Range('E14').Value = '  -1.46%  '
$ws.Range('E15').Value = '  -3.23%  '
$ws.Range('D16').Value = "'2.175.51"
$ws.Range('E16').Value = '  -1.55%  '
$ws.Range('D17').Value = "'0.771"
$ws.Range('E17').Value = '  -3.33%  '
$ws.Range('D18').Value = "'42.367.29"
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('E19').Value = '  -2.90%  '
$ws.Range('D20').Value = "'70.62"
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('D21').Value = "'5.88"
$ws.Range('E21').Value = '  -0.67%  '
$ws.Range('D22').Value = "'225.33"
$ws.Range('E22').Value = '  -1.45%  '
$ws.Range('D23').Value = "'9.41"
$ws.Range('E23').Value = '  -6.26%  '
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('D26').Value = "'10.46"
$ws.Range('E26').Value = '  -4.12%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('E28').Value = '  +1.65%  '
$ws.Range('E29').Value = '  -1.98%  '
$ws.Range('D30').Value = "'170.78"
$ws.Range('E30').Value = '  -1.15%  '
$ws.Range('D31').Value = "'36.75"
$ws.Range('E31').Value = '  +6.84%  '
$ws.Range('D32').Value = "'20.02"
$ws.Range('E32').Value = '  -0.95%  '
$ws.Range('E33').Value = '  +1.65%  '
$ws.Range('E34').Value = '  -4.19%  '
$ws.Range('E35').Value = '  -1.07%  '
$ws.Range('D36').Value = "'0.108"
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('D37').Value = "'4.25"
$ws.Range('E37').Value = '  -3.40%  '
$ws.Range('D38').Value = "'0.0339"
$ws.Range('E38').Value = '  +5.33%  '
$ws.Range('D39').Value = "'11.85"
$ws.Range('E39').Value = '  -5.98%  '
$ws.Range('E40').Value = '  -3.39%  '
$ws.Range('E41').Value = '  -0.66%  '
$ws.Range('D42').Value = "'59.38"
$ws.Range('E42').Value = '  -2.04%  '
$ws.Range('E43').Value = '  -6.30%  '
$ws.Range('D44').Value = "'102.30"
$ws.Range('E44').Value = '  +2.75%  '
$ws.Range('D45').Value = "'2.46"
$ws.Range('E45').Value = '  +6.54%  '
$ws.Range('E46').Value = '  +9.81%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = "'0.0969"
$ws.Range('E47').Value = '  -1.15%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').Value = "'8.24"
$ws.Range('E48').Value = '  -3.66%  '
$ws.Range('E49').Value = '  -1.18%  '
$ws.Range('D50').Value = "'1.13"
$ws.Range('E50').Value = '  -0.98%  '
$ws.Range('E51').Value = '  +0.42%  '
